$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 (P03.01 / čerpadlo nádrží) will become the new last row of the table
# once the rows below it are removed, so give it in advance the "closing"
# bottom-border formatting that currently belongs to row 29 (the current
# last row of the table).
$ws.Range("B29:C29").Copy()
$ws.Range("B26:C26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Remove the A01.01 / "pasterizační tank" row - everything below shifts up
$ws.Rows.Item(7).Delete()

# Remove the three tank rows (T03.01, T03.02, T03.03), which after the shift
# above now sit at rows 26-28, leaving P03.01 as the new, final row (25).
$ws.Range("A26:A28").EntireRow.Delete()

# Match the author's final cell selection
$ws.Range("E13").Select()
